$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column from 2023-09-14 to 2023-09-15
# for rows 2 through 5 (cells C2:C5), keeping existing formatting.
$newDate = Get-Date -Year 2023 -Month 9 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
